$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 26
    $ws.Range("B26").Value = 5143628
    $ws.Range("F26").Value = 'FC Nordsjaelland'
    $ws.Range("G26").Value = 'Brondby'
    $ws.Range("K26").Value = 2.2
    $ws.Range("L26").Value = 3.5
    $ws.Range("M26").Value = 3.1
    $ws.Range("N26").Value = 2.15
    $ws.Range("O26").Value = 3.5
    $ws.Range("P26").Value = 3.2
    $ws.Range("Q26").Value = -0.25
    $ws.Range("R26").Value = 1.875
    $ws.Range("S26").Value = 1.975
    $ws.Range("W26").Value = 1.15
    $ws.Range("Z26").Value = 0.875
    $ws.Range("AA26").Value = -1
    # Row 27
    $ws.Range("B27").Value = 5143627
    $ws.Range("F27").Value = 'FC Copenhagen'
    $ws.Range("G27").Value = 'Viborg'
    $ws.Range("K27").Value = 1.444
    $ws.Range("L27").Value = 4.333
    $ws.Range("M27").Value = 6.5
    $ws.Range("N27").Value = 1.333
    $ws.Range("O27").Value = 4.75
    $ws.Range("P27").Value = 9.5
    $ws.Range("Q27").Value = -1.5
    $ws.Range("R27").Value = 2.025
    $ws.Range("S27").Value = 1.825
    $ws.Range("W27").Value = 0.333
    $ws.Range("Z27").Value = -1
    $ws.Range("AA27").Value = 0.825
    # Row 33
    $ws.Range("B33").Value = 6433595
    $ws.Range("F33").Value = 'Lyngby'
    $ws.Range("G33").Value = 'Silkeborg IF'
    $ws.Range("H33").Value = 1
    $ws.Range("I33").Value = 1
    $ws.Range("J33").Value = 'D'
    $ws.Range("K33").Value = 3.1
    $ws.Range("L33").Value = 3.6
    $ws.Range("M33").Value = 2.15
    $ws.Range("N33").Value = 2.5
    $ws.Range("O33").Value = 3.6
    $ws.Range("P33").Value = 2.7
    $ws.Range("R33").Value = 1.825
    $ws.Range("S33").Value = 2.025
    $ws.Range("T33").Value = 2.5
    $ws.Range("U33").Value = 1.825
    $ws.Range("V33").Value = 2.025
    $ws.Range("X33").Value = 2.6
    $ws.Range("Y33").Value = -1
    $ws.Range("Z33").Value = 0
    $ws.Range("AA33").Value = -0
    $ws.Range("AB33").Value = -1
    $ws.Range("AC33").Value = 1.025
    # Row 34
    $ws.Range("B34").Value = 6433596
    $ws.Range("F34").Value = 'AC Horsens'
    $ws.Range("G34").Value = 'AaB'
    $ws.Range("H34").Value = 0
    $ws.Range("I34").Value = 4
    $ws.Range("J34").Value = 'A'
    $ws.Range("K34").Value = 2.45
    $ws.Range("L34").Value = 3.4
    $ws.Range("M34").Value = 2.75
    $ws.Range("N34").Value = 2.625
    $ws.Range("O34").Value = 3.25
    $ws.Range("P34").Value = 2.75
    $ws.Range("R34").Value = 1.875
    $ws.Range("S34").Value = 1.975
    $ws.Range("T34").Value = 2.25
    $ws.Range("U34").Value = 1.9
    $ws.Range("V34").Value = 1.95
    $ws.Range("X34").Value = -1
    $ws.Range("Y34").Value = 1.75
    $ws.Range("Z34").Value = -1
    $ws.Range("AA34").Value = 0.9750000000000001
    $ws.Range("AB34").Value = 0.8999999999999999
    $ws.Range("AC34").Value = -1
    # Row 45
    $ws.Range("B45").Value = 6478386
    $ws.Range("F45").Value = 'Lyngby'
    $ws.Range("G45").Value = 'AC Horsens'
    $ws.Range("I45").Value = 1
    $ws.Range("K45").Value = 2.05
    $ws.Range("L45").Value = 3.4
    $ws.Range("M45").Value = 3.5
    $ws.Range("N45").Value = 2.15
    $ws.Range("O45").Value = 3.4
    $ws.Range("P45").Value = 3.4
    $ws.Range("Q45").Value = -0.25
    $ws.Range("R45").Value = 1.875
    $ws.Range("S45").Value = 1.975
    $ws.Range("T45").Value = 2.25
    $ws.Range("U45").Value = 1.875
    $ws.Range("V45").Value = 1.975
    $ws.Range("W45").Value = 1.15
    $ws.Range("Z45").Value = 0.875
    $ws.Range("AB45").Value = 0.875
    $ws.Range("AC45").Value = -1
    # Row 46
    $ws.Range("B46").Value = 6445249
    $ws.Range("F46").Value = 'Odense BK'
    $ws.Range("G46").Value = 'Silkeborg IF'
    $ws.Range("I46").Value = 0
    $ws.Range("K46").Value = 2.75
    $ws.Range("L46").Value = 3.5
    $ws.Range("M46").Value = 2.45
    $ws.Range("N46").Value = 2.45
    $ws.Range("O46").Value = 3.6
    $ws.Range("P46").Value = 2.7
    $ws.Range("Q46").Value = 0
    $ws.Range("R46").Value = 1.825
    $ws.Range("S46").Value = 2.025
    $ws.Range("T46").Value = 2.75
    $ws.Range("U46").Value = 1.85
    $ws.Range("V46").Value = 2
    $ws.Range("W46").Value = 1.45
    $ws.Range("Z46").Value = 0.825
    $ws.Range("AB46").Value = -1
    $ws.Range("AC46").Value = 1
    # Row 69
    $ws.Range("B69").Value = 6471205
    $ws.Range("F69").Value = 'Midtjylland'
    $ws.Range("G69").Value = 'AC Horsens'
    $ws.Range("H69").Value = 3
    $ws.Range("I69").Value = 1
    $ws.Range("K69").Value = 1.416
    $ws.Range("L69").Value = 4.333
    $ws.Range("M69").Value = 6.5
    $ws.Range("N69").Value = 1.363
    $ws.Range("O69").Value = 5
    $ws.Range("P69").Value = 9
    $ws.Range("Q69").Value = -1.5
    $ws.Range("R69").Value = 1.95
    $ws.Range("S69").Value = 1.9
    $ws.Range("T69").Value = 3
    $ws.Range("U69").Value = 2
    $ws.Range("V69").Value = 1.85
    $ws.Range("W69").Value = 0.363
    $ws.Range("Z69").Value = 0.95
    $ws.Range("AB69").Value = 1
    $ws.Range("AC69").Value = -1
    # Row 70
    $ws.Range("B70").Value = 6445253
    $ws.Range("F70").Value = 'Silkeborg IF'
    $ws.Range("G70").Value = 'Lyngby'
    $ws.Range("H70").Value = 1
    $ws.Range("I70").Value = 0
    $ws.Range("K70").Value = 1.85
    $ws.Range("L70").Value = 3.6
    $ws.Range("M70").Value = 3.75
    $ws.Range("N70").Value = 1.909
    $ws.Range("O70").Value = 3.6
    $ws.Range("P70").Value = 4.2
    $ws.Range("Q70").Value = -0.5
    $ws.Range("R70").Value = 1.9
    $ws.Range("S70").Value = 1.95
    $ws.Range("T70").Value = 2.75
    $ws.Range("U70").Value = 1.975
    $ws.Range("V70").Value = 1.875
    $ws.Range("W70").Value = 0.909
    $ws.Range("Z70").Value = 0.8999999999999999
    $ws.Range("AB70").Value = -1
    $ws.Range("AC70").Value = 0.875
    # Row 130
    $ws.Range("B130").Value = 6779623
    $ws.Range("F130").Value = 'Silkeborg IF'
    $ws.Range("G130").Value = 'Hvidovre IF'
    $ws.Range("I130").Value = 0
    $ws.Range("J130").Value = 'H'
    $ws.Range("K130").Value = 1.533
    $ws.Range("L130").Value = 4
    $ws.Range("M130").Value = 6
    $ws.Range("N130").Value = 1.444
    $ws.Range("O130").Value = 4.75
    $ws.Range("P130").Value = 7
    $ws.Range("Q130").Value = -1.25
    $ws.Range("R130").Value = 1.875
    $ws.Range("S130").Value = 1.975
    $ws.Range("T130").Value = 3
    $ws.Range("U130").Value = 1.9
    $ws.Range("V130").Value = 1.95
    $ws.Range("W130").Value = 0.444
    $ws.Range("X130").Value = -1
    $ws.Range("Z130").Value = -0.5
    $ws.Range("AA130").Value = 0.4875
    $ws.Range("AC130").Value = 0.95
    # Row 131
    $ws.Range("B131").Value = 6779624
    $ws.Range("F131").Value = 'Lyngby'
    $ws.Range("G131").Value = 'FC Nordsjaelland'
    $ws.Range("I131").Value = 1
    $ws.Range("J131").Value = 'D'
    $ws.Range("K131").Value = 4.5
    $ws.Range("L131").Value = 3.6
    $ws.Range("M131").Value = 1.75
    $ws.Range("N131").Value = 4.5
    $ws.Range("O131").Value = 3.6
    $ws.Range("P131").Value = 1.8
    $ws.Range("Q131").Value = 0.75
    $ws.Range("R131").Value = 1.825
    $ws.Range("S131").Value = 2.025
    $ws.Range("T131").Value = 2.75
    $ws.Range("U131").Value = 2
    $ws.Range("V131").Value = 1.85
    $ws.Range("W131").Value = -1
    $ws.Range("X131").Value = 2.6
    $ws.Range("Z131").Value = 0.825
    $ws.Range("AA131").Value = -1
    $ws.Range("AC131").Value = 0.8500000000000001
    # Row 196
    $ws.Range("N196").Value = 1.45
    $ws.Range("P196").Value = 5.75
    $ws.Range("Q196").Value = -1.25
    $ws.Range("R196").Value = 2.06
    $ws.Range("S196").Value = 1.84
    $ws.Range("T196").Value = 3
    $ws.Range("U196").Value = 2.05
    $ws.Range("V196").Value = 1.8
    # Row 197
    $ws.Range("R197").Value = 1.89
    $ws.Range("S197").Value = 2.01
    # Row 199
    $ws.Range("R199").Value = 2.05
    $ws.Range("S199").Value = 1.85
